$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 328.86957
$ws.Range("I2").Value = 206.33333
$ws.Range("K2").Value = 206.33333
$ws.Range("M2").Value = -93.33332999999999
$ws.Range("H8").Value = 1531.1428
$ws.Range("I8").Value = 119.833336
$ws.Range("K8").Value = 359.500008
$ws.Range("M8").Value = -220.500008
$ws.Range("H9").Value = 339.3
$ws.Range("I9").Value = 354.77777
$ws.Range("J9").Value = 200
$ws.Range("K9").Value = 354.77777
$ws.Range("L9").Value = 200
$ws.Range("M9").Value = -185.77777
$ws.Range("N9").Value = -538
$ws.Range("H21").Value = 9998
$ws.Range("I21").Value = 9998
$ws.Range("K21").Value = 9998
$ws.Range("M21").Value = -9530
$ws.Range("H23").Value = 9998
$ws.Range("I23").Value = 9998
$ws.Range("K23").Value = 9998
$ws.Range("M23").Value = -9764
$ws.Range("H92").Value = 966.3333
$ws.Range("I92").Value = 919.6
$ws.Range("J92").Value = 1200
$ws.Range("K92").Value = 919.6
$ws.Range("L92").Value = 1200
$ws.Range("M92").Value = 328.4
$ws.Range("N92").Value = -3696
$ws.Range("H129").Value = 1298
$ws.Range("I129").Value = 1298
$ws.Range("K129").Value = 3894
$ws.Range("M129").Value = 1106
$ws.Range("H138").Value = 8540.154
$ws.Range("I138").Value = 8947.333000000001
$ws.Range("K138").Value = 26841.999
$ws.Range("M138").Value = -21701.999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3478.889
$ws.Range("I45").Value = 2902.75
$ws.Range("J45").Value = 3939.8
$ws.Range("K45").Value = 2902.75
$ws.Range("L45").Value = 3939.8
$ws.Range("M45").Value = -2525.75
$ws.Range("N45").Value = -4693.8
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()  # was -34630
$ws.Range("H61").Value = 4003640.8
$ws.Range("I61").Value = 3690.3333
$ws.Range("K61").Value = 3690.3333
$ws.Range("M61").Value = -3478.3333
$ws.Range("H74").Value = 3099874.2
$ws.Range("I74").Value = 4280537.5
$ws.Range("J74").Value = 30150
$ws.Range("K74").Value = 4280537.5
$ws.Range("L74").Value = 30150
$ws.Range("M74").Value = -4279663.5
$ws.Range("N74").Value = -31898
$ws.Range("H77").Value = 3099874.2
$ws.Range("I77").Value = 4280537.5
$ws.Range("J77").Value = 30150
$ws.Range("K77").Value = 21402687.5
$ws.Range("L77").Value = 150750
$ws.Range("M77").Value = -21398319.5
$ws.Range("N77").Value = -159486
$ws.Range("H102").Value = 4952.857
$ws.Range("I102").Value = 3421.3076
$ws.Range("K102").Value = 3421.3076
$ws.Range("M102").Value = -1799.3076
$ws.Range("H122").Value = 1438.5454
$ws.Range("I122").Value = 1234.5
$ws.Range("K122").Value = 3703.5
$ws.Range("M122").Value = -1253.5
$ws.Range("H132").Value = 611737.0600000001
$ws.Range("I132").Value = 659684.75
$ws.Range("K132").Value = 1979054.25
$ws.Range("M132").Value = -1976524.25
$ws.Range("H136").Value = 4003640.8
$ws.Range("I136").Value = 3690.3333
$ws.Range("K136").Value = 11070.9999
$ws.Range("M136").Value = -8520.999899999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 1990.8334
$ws.Range("I5").Value = 1389
$ws.Range("K5").Value = 1389
$ws.Range("M5").Value = -1276
$ws.Range("H94").Value = 3471.7856
$ws.Range("I94").Value = 3334.889
$ws.Range("J94").Value = 3718.2
$ws.Range("K94").Value = 3334.889
$ws.Range("L94").Value = 3718.2
$ws.Range("M94").Value = -2883.889
$ws.Range("N94").Value = -4620.2
$ws.Range("H106").Value = 14780.667
$ws.Range("J106").Value = 14780.667
$ws.Range("L106").Value = 14780.667
$ws.Range("N106").Value = -17304.667
$ws.Range("H134").Value = 5753831.5
$ws.Range("J134").Value = 41696668
$ws.Range("L134").Value = 125090004
$ws.Range("N134").Value = -125095074

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 597178.1
$ws.Range("I31").Value = 869836.7
$ws.Range("J31").Value = 2286.682
$ws.Range("K31").Value = 869836.7
$ws.Range("L31").Value = 2286.682
$ws.Range("M31").Value = -869541.7
$ws.Range("N31").Value = -2876.682
$ws.Range("H34").Value = 597178.1
$ws.Range("I34").Value = 869836.7
$ws.Range("J34").Value = 2286.682
$ws.Range("K34").Value = 869836.7
$ws.Range("L34").Value = 2286.682
$ws.Range("M34").Value = -869634.7
$ws.Range("N34").Value = -2690.682
$ws.Range("H58").Value = 7523109
$ws.Range("I58").Value = 7250226
$ws.Range("K58").Value = 7250226
$ws.Range("M58").Value = -7250023
$ws.Range("H132").Value = 1483.36
$ws.Range("I132").Value = 1483.36
$ws.Range("K132").Value = 4450.08
$ws.Range("M132").Value = -1920.08
$ws.Range("H134").Value = 1651.9429
$ws.Range("I134").Value = 1412.1731
$ws.Range("K134").Value = 4236.5193
$ws.Range("M134").Value = -1701.5193
$ws.Range("H136").Value = 7523109
$ws.Range("I136").Value = 7250226
$ws.Range("K136").Value = 21750678
$ws.Range("M136").Value = -21748128

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2809583.2
$ws.Range("I5").Value = 2748304.2
$ws.Range("J5").Value = 2909161.5
$ws.Range("K5").Value = 8244912.600000001
$ws.Range("L5").Value = 8727484.5
$ws.Range("M5").Value = -8244800.600000001
$ws.Range("N5").Value = -8727708.5
$ws.Range("H37").Value = 175000
$ws.Range("J37").Value = 175000
$ws.Range("L37").Value = 525000
$ws.Range("N37").Value = -525224
$ws.Range("H57").Value = 600
$ws.Range("I57").Value = 600
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 1800
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -1241
$ws.Range("N57").ClearContents()  # was -4718
$ws.Range("H113").Value = 1142.1818
$ws.Range("I113").Value = 426
$ws.Range("J113").Value = 1410.75
$ws.Range("K113").Value = 1278
$ws.Range("L113").Value = 4232.25
$ws.Range("M113").Value = 892
$ws.Range("N113").Value = -8572.25
$ws.Range("H135").Value = 2809583.2
$ws.Range("I135").Value = 2748304.2
$ws.Range("J135").Value = 2909161.5
$ws.Range("K135").Value = 24734737.8
$ws.Range("L135").Value = 26182453.5
$ws.Range("M135").Value = -24732202.8
$ws.Range("N135").Value = -26187523.5
$ws.Range("H140").Value = 1442.7142
$ws.Range("I140").Value = 1442.7142
$ws.Range("K140").Value = 4328.142599999999
$ws.Range("M140").Value = 851.8574000000008

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 130.22223
$ws.Range("I2").Value = 60.357143
$ws.Range("K2").Value = 60.357143
$ws.Range("M2").Value = 52.642857
$ws.Range("H15").Value = 24766.666
$ws.Range("I15").Value = 25000
$ws.Range("K15").Value = 25000
$ws.Range("M15").Value = -24712
$ws.Range("H81").Value = 24766.666
$ws.Range("I81").Value = 25000
$ws.Range("K81").Value = 25000
$ws.Range("M81").Value = -24002
$ws.Range("H84").Value = 24766.666
$ws.Range("I84").Value = 25000
$ws.Range("K84").Value = 75000
$ws.Range("M84").Value = -70008
$ws.Range("H98").Value = 28357.2
$ws.Range("J98").Value = 28357.2
$ws.Range("L98").Value = 28357.2
$ws.Range("N98").Value = -34347.2
$ws.Range("H102").Value = 2563
$ws.Range("I102").Value = 2780.1333
$ws.Range("J102").Value = 1748.75
$ws.Range("K102").Value = 2780.1333
$ws.Range("L102").Value = 1748.75
$ws.Range("M102").Value = -1158.1333
$ws.Range("N102").Value = -4992.75
$ws.Range("H113").Value = 3349.5
$ws.Range("I113").Value = 3166.111
$ws.Range("K113").Value = 3166.111
$ws.Range("M113").Value = -996.1109999999999
$ws.Range("H122").Value = 8389.5
$ws.Range("I122").Value = 9626
$ws.Range("K122").Value = 28878
$ws.Range("M122").Value = -26428
$ws.Range("H132").Value = 13832.707
$ws.Range("I132").Value = 14490.615
$ws.Range("K132").Value = 43471.845
$ws.Range("M132").Value = -40941.845

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2362.7222
$ws.Range("J16").Value = 6742.6
$ws.Range("L16").Value = 6742.6
$ws.Range("N16").Value = -7082.6
$ws.Range("H136").Value = 20835784
$ws.Range("I136").Value = 11366309
$ws.Range("K136").Value = 34098927
$ws.Range("M136").Value = -34096377

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 17005000
$ws.Range("I49").Value = 25500000
$ws.Range("K49").Value = 25500000
$ws.Range("M49").Value = -25499770
$ws.Range("H97").Value = 27333.334
$ws.Range("J97").Value = 27333.334
$ws.Range("L97").Value = 27333.334
$ws.Range("N97").Value = -29315.334
$ws.Range("H101").Value = 9010
$ws.Range("J101").Value = 9010
$ws.Range("L101").Value = 9010
$ws.Range("N101").Value = -15500
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()  # was -41988
$ws.Range("H132").Value = 7578301.5
$ws.Range("I132").Value = 11113423
$ws.Range("J132").Value = 3041.1428
$ws.Range("K132").Value = 33340269
$ws.Range("L132").Value = 9123.428400000001
$ws.Range("M132").Value = -33337739
$ws.Range("N132").Value = -14183.4284
